$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 16 (this pushes the old row 16 and
# everything below it down by two rows, matching the target diff).
$ws.Rows.Item(16).Resize(2).Insert()

# --- New row 16: "US Citizen Indicator" ---------------------------------
$ws.Cells.Item(16, 3).Value = "US Citizen Indicator"
$ws.Cells.Item(16, 5).Value = $true
$ws.Cells.Item(16, 6).Value = "wm-req-doc:WarrantModificationRequest/nc:Person/nc:PersonUSCitizenIndicator"

# --- New row 17: "PersonImmigationAlienQueryIndicator" ------------------
$ws.Cells.Item(17, 1).Value = "ext"
$ws.Cells.Item(17, 3).Value = "PersonImmigationAlienQueryIndicator"
$ws.Cells.Item(17, 4).Value = "True if a person query should include the DHS ICE Database; false otherwise."
$ws.Cells.Item(17, 5).Value = $false
$ws.Cells.Item(17, 6).Value = "wm-req-doc:WarrantModificationRequest/nc:Person/wm-req-ext:PersonImmigrationAlienQueryIndicator"

# Match styling used by the neighboring rows for the new content.
$ws.Cells.Item(16, 2).Style = $ws.Cells.Item(15, 2).Style
$ws.Cells.Item(16, 3).Style = $ws.Cells.Item(8, 3).Style
$ws.Cells.Item(16, 5).Style = $ws.Cells.Item(15, 2).Style
$ws.Cells.Item(16, 6).Style = $ws.Cells.Item(15, 2).Style

$ws.Cells.Item(17, 1).Style = $ws.Cells.Item(8, 1).Style
$ws.Cells.Item(17, 2).Style = $ws.Cells.Item(15, 2).Style
$ws.Cells.Item(17, 3).Style = $ws.Cells.Item(8, 3).Style
$ws.Cells.Item(17, 4).Style = $ws.Cells.Item(8, 3).Style
$ws.Cells.Item(17, 5).Style = $ws.Cells.Item(15, 2).Style
$ws.Cells.Item(17, 6).Style = $ws.Cells.Item(15, 2).Style

$ws.Rows.Item(17).RowHeight = 30

# Restore the sheet view state (frozen pane position and selection) to
# match what results from scrolling/selecting rows 16:17 before inserting.
$ws.Range("A16:XFD17").Select()
$excel.ActiveWindow.ScrollRow = 14
